$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update "VALOR MORA" total (E11): 170820 -> 56940
$ws.Range("E11").Value = 56940

# Update "Cant. Trabajadores" (C13): 3 -> 1
$ws.Range("C13").Value = 1

# Update the single remaining worker row (row 16) to the KEYLA record,
# which is what survives after removing the other two old EC rows.
$ws.Range("C16").Value = "1052987601"
$ws.Range("D16").Value = "KEYLA XIMENA TORRES MORALES"
$ws.Range("E16").Value = "2509"

# Remove the two now-obsolete worker rows (17 and 18); this shifts the
# trailing signature rows (23/24) up to (21/22).
$ws.Rows("17:18").Delete()
